$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New dates on row 2 (H2, I2), copying date-format style from G2 ---
$ws.Range("H2").Value = 43993
$ws.Range("I2").Value = 43996
$ws.Range("G2").Copy()
$ws.Range("H2:I2").PasteSpecial(-4122)  # xlPasteFormats

# --- New row 12: "Algoritmo Q-learning" task with hour entries ---
# Shared-string order must match: "3.5 h." (H12), "Algoritmo Q-learning" (A12), "4.5 h." (I12)
$ws.Range("H12").Value = "3.5 h."
$ws.Range("A12").Value = "Algoritmo Q-learning"
$ws.Range("I12").Value = "4.5 h."

$ws.Range("A9").Copy()
$ws.Range("A12").PasteSpecial(-4122)  # xlPasteFormats, matches other task rows (style 3)

$ws.Range("G11").Copy()
$ws.Range("H12:I12").PasteSpecial(-4122)  # xlPasteFormats, matches hour entry style (style 1)

# --- New total-hours label in A1 ---
$ws.Range("A1").Value = "Total horas: 29.5"

$ws.Application.CutCopyMode = $false
